# Edit script: update Percent-Change / Weight values and refresh the
# "as of" date in the confidential-use disclaimer, matching the
# upstream holdings-file diff. The worksheet ships with sheet
# protection enabled, so we briefly unprotect, write all the cell
# changes, then re-protect before returning control.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Update the disclaimer date (shared string used by cell A59) ---
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# --- Refresh Weight (D) / Percent Change (E) figures for each holding row ---
$ws.Range("D2").Value = 0.01430595795285261
$ws.Range("E2").Value = -0.00235183443085607
$ws.Range("D3").Value = 0.05094123661787659
$ws.Range("E3").Value = -0.01165304443812498
$ws.Range("D4").Value = 0.01410981863873646
$ws.Range("E4").Value = 0.003731343283582156
$ws.Range("D5").Value = 0.009810704066305057
$ws.Range("E5").Value = -0.01257462212625438
$ws.Range("D6").Value = 0.01548235769549141
$ws.Range("E6").Value = -0.004515290415269768
$ws.Range("D7").Value = 0.02016630532816372
$ws.Range("E7").Value = -0.01023742104116743
$ws.Range("D8").Value = 0.004908342721550742
$ws.Range("E8").Value = -0.02838990581126721
$ws.Range("D9").Value = 0.006974659180682845
$ws.Range("E9").Value = -0.02830036983437856
$ws.Range("D10").Value = 0.0144957732071267
$ws.Range("E10").Value = -0.007794820216243448
$ws.Range("D11").Value = 0.00832595188860882
$ws.Range("E11").Value = -0.005163511187607495
$ws.Range("D12").Value = 0.01603756653541802
$ws.Range("E12").Value = -0.02074592074592074
$ws.Range("D13").Value = 0.003000657359540528
$ws.Range("E13").Value = 0.01744186046511631
$ws.Range("D14").Value = 0.005955831333024898
$ws.Range("E14").Value = -0.001726121979286344
$ws.Range("D15").Value = 0.01487689905987398
$ws.Range("E15").Value = -0.0140887836278617
$ws.Range("D16").Value = 0.01085339511514041
$ws.Range("E16").Value = -0.01033324722293982
$ws.Range("D17").Value = 0.02102827779993622
$ws.Range("E17").Value = 0.003599999999999826
$ws.Range("D18").Value = 0.008564957736432685
$ws.Range("E18").Value = 0.00213871066300042
$ws.Range("D19").Value = 0.01698611320571796
$ws.Range("E19").Value = 0.0003521333411584671
$ws.Range("D20").Value = 0.01218269380213786
$ws.Range("E20").Value = -0.004602874239247212
$ws.Range("D21").Value = 0.00723179607692088
$ws.Range("E21").Value = -0.05801721389862924
$ws.Range("D22").Value = 0.01475206897225274
$ws.Range("E22").Value = 0.04789498580889284
$ws.Range("D23").Value = 0.01998106956549939
$ws.Range("E23").Value = -0.02045885076825937
$ws.Range("D24").Value = 0.01043289186514998
$ws.Range("E24").Value = -0.01887772774506413
$ws.Range("D25").Value = 0.02012309611141037
$ws.Range("E25").Value = -0.012719368861697
$ws.Range("D26").Value = 0.0141499748610833
$ws.Range("E26").Value = -0.01063387119145598
$ws.Range("D27").Value = 0.02028923508253194
$ws.Range("E27").Value = -0.01222828381779184
$ws.Range("D28").Value = 0.05507165767284629
$ws.Range("E28").Value = -0.01124574324859429
$ws.Range("D29").Value = 0.02078743392486761
$ws.Range("E29").Value = -0.01832993890020362
$ws.Range("D30").Value = 0.0287529767084876
$ws.Range("E30").Value = -0.005504030510531499
$ws.Range("D31").Value = 0.01496154177214425
$ws.Range("E31").Value = -0.006856702619414512
$ws.Range("D32").Value = 0.01313796952133378
$ws.Range("E32").Value = 0.0008109569291765428
$ws.Range("D33").Value = 0.01803160802492664
$ws.Range("E33").Value = -0.01878343077497202
$ws.Range("D34").Value = 0.04278404055273778
$ws.Range("E34").Value = -0.01155566817538412
$ws.Range("D35").Value = 0.01090105921148693
$ws.Range("E35").Value = 0.005829903978052053
$ws.Range("D36").Value = 0.01019612903362478
$ws.Range("E36").Value = 0.002737616562580092
$ws.Range("D37").Value = 0.01092302207941131
$ws.Range("E37").Value = -0.0004278074866311821
$ws.Range("D38").Value = 0.007533730993123816
$ws.Range("E38").Value = -0.009304056568663954
$ws.Range("D39").Value = 0.01238275839477666
$ws.Range("E39").Value = -0.01357042583060375
$ws.Range("D40").Value = 0.01748194441973838
$ws.Range("E40").Value = -0.0001888930865129312
$ws.Range("D41").Value = 0.01745297212588069
$ws.Range("E41").Value = -0.02102692803364292
$ws.Range("D42").Value = 0.031795690623628
$ws.Range("E42").Value = 0.0004702977572677014
$ws.Range("D43").Value = 0.01144156383345567
$ws.Range("E43").Value = -0.008070356958096214
$ws.Range("D44").Value = 0.02186828740378434
$ws.Range("E44").Value = -0.003842077371489228
$ws.Range("D45").Value = 0.01233490738040525
$ws.Range("E45").Value = -0.01409283084056923
$ws.Range("D46").Value = 0.008670940256544366
$ws.Range("E46").Value = -0.01037961585720748
$ws.Range("D47").Value = 0.0137266989937237
$ws.Range("E47").Value = -0.02173291392622256
$ws.Range("D48").Value = 0.01092545201373486
$ws.Range("E48").Value = -0.02194164292863177
$ws.Range("D49").Value = 0.0158779385421634
$ws.Range("E49").Value = -0.01097167612365502
$ws.Range("D50").Value = 0.00865791830081048
$ws.Range("E50").Value = -0.00666748226082714
$ws.Range("D51").Value = 0.01132318241769306
$ws.Range("E51").Value = -0.0356508102456875
$ws.Range("D52").Value = 0.008261745547059681
$ws.Range("E52").Value = 0.004521133186776716
$ws.Range("D53").Value = 0.009823227573972573
$ws.Range("E53").Value = -0.002219952936997727
$ws.Range("D54").Value = 0.1352815466768732
$ws.Range("E54").Value = 0.0004928050463237632
$ws.Range("D55").Value = 0.04364442422129872
$ws.Range("E55").Value = -0.006431258342434187
$ws.Range("E56").Value = -0.00741145545184263

$ws.Protect()

